$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the Joel Embiid / 2023 row (row 7) entirely - this shifts all
# subsequent rows up by one.
$ws.Rows.Item(7).Delete()

# Append the new Oscar Robertson / 1964 record in the now-empty last row (11).
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "1964"
$ws.Range("A11").Style = "Normal"
$ws.Range("B11").Value = "nba mvp"
$ws.Range("C11").Value = "Oscar Robertson"
$ws.Range("D11").Value = 25
$ws.Range("E11").Value = "CIN"
$ws.Range("F11").Value = 3.5
$ws.Range("G11").Value = 31.4
$ws.Range("H11").Value = "1963-64"
$ws.Range("I11").Value = "Yes"
$ws.Range("J11").Value = "No"
$ws.Range("K11").Value = 1964
